# Update loading_percent values for columns B, C, E, F, G, J, K, L across rows 2-25
# (case with 380 kV done)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters (in order) whose values change, mapped to column indices:
# B=2, C=3, E=5, F=6, G=7, J=10, K=11, L=12
$cols = @(2,3,5,6,7,10,11,12)

# Each row below holds the new values for columns B,C,E,F,G,J,K,L respectively,
# for worksheet rows 2 through 25 (in order).
$data = @(
    @(19.6276857286398,4.378135040345265,9.527082579122737,63.66074658758853,3.808216604563695,13.23575416922598,16.39020577852813,9.658448768736228),
    @(19.64289180060902,4.154497746155063,9.593568983954647,62.74380268717346,3.811867872428851,13.12492447764492,16.39070802306887,9.747725960464383),
    @(19.65911656871898,4.010131143323904,9.637186216941181,62.17729262793009,3.814224011442,13.05558538531232,16.39610567403536,9.805351834431047),
    @(19.66745736894361,3.949557903537658,9.6556623688823,61.94573438778102,3.815212999188835,13.02701201964488,16.39958427808953,9.829544300636321),
    @(19.66894667267013,3.939395531159591,9.658772680924523,61.90724727127409,3.815378965177891,13.02224848641047,16.40023910791102,9.83360437939842),
    @(19.65922205968581,4.009321242252594,9.63743255248885,62.17417235173713,3.814237232347623,13.05520130954516,16.39614741072797,9.805675225983444),
    @(19.63149768951435,4.302505179689088,9.549426589482289,63.34541890312283,3.80945191876569,13.19781189568867,16.38932251437045,9.688650297579361),
    @(19.63186532613679,4.820422973466304,9.399060133361917,65.60663268066054,3.800969219412329,13.46693700925602,16.41631216362402,9.481317479484261),
    @(19.6655333703286,5.165103295453145,9.302187641097161,67.23601401097842,3.795279137910111,13.65784337454304,16.4606850580003,9.342302137310316),
    @(19.68808091465969,5.313977128199861,9.26108591880344,67.96823683775784,3.792806736634379,13.74312723301751,16.48616574266204,9.281909790005262),
    @(19.69765499463402,5.36920371831941,9.245949703964179,68.24404888227681,3.79188706959654,13.77519044755929,16.49657144091855,9.259446895803833),
    @(19.69554705698031,5.357360896267211,9.249190491648655,68.18471564557824,3.792084400792248,13.76829549498865,16.49429681837564,9.264266654916348),
    @(19.68884784910206,5.318543719702393,9.25983206374519,67.99095838154913,3.792730743459069,13.74576978410549,16.48700669064821,9.280053626576171),
    @(19.68487913922442,5.294617247209209,9.266406136204841,67.87208069699646,3.793128802661481,13.73194170910407,16.48263966623682,9.289776437813229),
    @(19.66420503733618,5.155213750129187,9.304933554679694,67.18796576894087,3.79544304090858,13.65223782715024,16.45912600722216,9.346305945524778),
    @(19.65337189015269,5.067657465508748,9.32932966646058,66.76585638095338,3.796892395088126,13.60293702084276,16.44605432545878,9.381711914920116),
    @(19.64782197304162,5.016552043521306,9.343640743156943,66.52223388223604,3.797736954545504,13.57443407446567,16.43903459935455,9.402344586594126),
    @(19.64605993541402,4.999120936554288,9.348534112524559,66.4396091855792,3.798024788318097,13.56475854922062,16.43674363089375,9.409376589721955),
    @(19.6544546395411,5.07705519946107,9.326703763972461,66.81087846938823,3.796736978523118,13.60820038524678,16.44739423985032,9.3779151637035),
    @(19.69078749663807,5.32997650342589,9.2566947439331,68.04791068228442,3.792540447907998,13.75239248564761,16.48912748346167,9.275405604722788),
    @(19.72056777133771,5.488576992555295,9.213435993828291,68.84776917215345,3.789894349280019,13.84527505618006,16.5208103970729,9.210777000731268),
    @(19.70412301599158,5.404544297658513,9.236295048845337,68.42171298219424,3.79129782124379,13.79582835740281,16.50349913275878,9.245054861423837),
    @(19.65396301568417,5.072808875878238,9.327890045508099,66.79052691715349,3.796807207124222,13.60582131375577,16.44678692119874,9.379630810879858),
    @(19.62589781197778,4.686543501789891,9.437354291069479,64.99979549056559,3.803168280362561,13.3953004357837,16.40469173463232,9.535054860341567)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Cells.Item($row, $cols[$j]).Value = $values[$j]
    }
}
